$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 7080
$ws.Range("D3").Value = 7081
$ws.Range("D4").Value = 7082
$ws.Range("D5").Value = 7083
$ws.Range("D6").Value = 7084
